$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 4.3
$ws.Range("M3").Value = 3.8
$ws.Range("M4").Value = 3.6
$ws.Range("M5").Value = 3.4
$ws.Range("M7").Value = 3.2
$ws.Range("M8").Value = 3
$ws.Range("M9").Value = 2.8
$ws.Range("M10").Value = 3.3
$ws.Range("M11").Value = 2.9
$ws.Range("M12").Value = 2.6
$ws.Range("M13").Value = 2.4
$ws.Range("M16").Value = 0.4
